# Actualización desde MV -datos-
# Appends six new daily rows (29-10-2021 .. 03-11-2021) to the bottom of the
# data table on Sheet1, continuing the existing B/C/D series with C reset to 0
# for the newly reported days. Mirrors the last row's "no D value" pattern for
# the very last appended row (03-11-2021), which has no D entry yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date labels for the new rows, in order (must remain plain text, not dates).
$newDates = @("29-10-2021", "30-10-2021", "31-10-2021", "01-11-2021", "02-11-2021", "03-11-2021")

$startRow = 303
$row = $startRow

foreach ($d in $newDates) {
    $cellA = $ws.Range("A" + $row)

    # Some of these labels (e.g. "01-11-2021") are still valid as a date in
    # either day/month order, so a plain .Value assignment would get silently
    # reinterpreted as a date serial. Force text format first, write the
    # value, then restore the default "Normal" style so the cell ends up
    # with no explicit style (matching the rest of the sheet).
    $cellA.NumberFormat = "@"
    $cellA.Value = $d
    $cellA.Style = "Normal"

    $ws.Range("B" + $row).Value = 17537
    $ws.Range("C" + $row).Value = 0

    # The final new row (03-11-2021) has not been reported for column D yet.
    if ($row -lt ($startRow + $newDates.Length - 1)) {
        $ws.Range("D" + $row).Value = 521
    }

    $row = $row + 1
}
